# The sheet previously had a two-line title split across A1 ("社員マスタ")
# and A2 ("2024年度版"). Replace it with a single consolidated title in A1
# ("社員一覧") and clear out the now-unused A2 cell, leaving the header row
# (名前/年齢/職業/出身地) and the data rows untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "社員一覧"
$ws.Range("A2").ClearContents()
